$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill column J (ReasonToReject) with "Nil" for rows 2 through 7
for ($r = 2; $r -le 7; $r++) {
    $ws.Cells.Item($r, 10).Value = "Nil"
}

# Update the active selection on the sheet to I7
$ws.Range("I7").Select()
